$wb = $excel.ActiveWorkbook

# --- ev_charging_uc: shuffle the two clustered timeslice-group strings ---
$ws = $wb.Worksheets.Item("ev_charging_uc")
$ws.Range("C13").Value = 'S1aH2,S2aH3,S3aH2,S3aH4,S4aH3,S4aH4,S4aH2,S1aH4,S2aH2,S1aH5,S2aH4,S2aH5,S5aH2,S5aH5,S5aH4,S1aH3,S3aH3,S4aH5,S3aH5,S5aH3'
$ws.Range("C14").Value = 'S3aH8,S2aH1,S3aH7,S5aH8,S4aH1,S1aH6,S2aH6,S3aH1,S1aH7,S3aH6,S4aH7,S5aH6,S5aH7,S4aH6,S1aH1,S2aH7,S4aH8,S2aH8,S1aH8,S5aH1'

# --- re_profiles: refreshed RE profile / ncap_afs sample values ---
$ws = $wb.Worksheets.Item("re_profiles")
$ws.Range("I4").Value = 0.069102744882046332
$ws.Range("N4").Value = 0.63189495603887091
$ws.Range("I5").Value = 0.0090283313245061378
$ws.Range("M5").Value = 'S1'
$ws.Range("N5").Value = 0.14770939379916706
$ws.Range("C6").Value = 0.000088889999999999998
$ws.Range("I6").Value = 0.009435445555022199
$ws.Range("M6").Value = 'S4'
$ws.Range("N6").Value = 0.18220731142989358
$ws.Range("C7").Value = 0.0095477800000000005
$ws.Range("I7").Value = 0.010138048170267647
$ws.Range("M7").Value = 'S5'
$ws.Range("N7").Value = 0.16266774641369736
$ws.Range("C8").Value = 0.11259713
$ws.Range("I8").Value = 0.086129717336578329
$ws.Range("M8").Value = 'S2'
$ws.Range("N8").Value = 0.075520592318371119
$ws.Range("C9").Value = 0.0026397600000000001
$ws.Range("I9").Value = 0.0099736640136369667
$ws.Range("C10").Value = 0
$ws.Range("I10").Value = 0.009999990864058169
$ws.Range("I11").Value = 0.040903710182698352
$ws.Range("C12").Value = 0
$ws.Range("I12").Value = 0.041396883574591439
$ws.Range("C13").Value = 0
$ws.Range("I13").Value = 0.0057188125975343321
$ws.Range("C14").Value = 0.0028831899999999999
$ws.Range("I14").Value = 0.0063792922521055265
$ws.Range("C15").Value = 0.0046609399999999997
$ws.Range("I15").Value = 0.0069000617624412229
$ws.Range("C16").Value = 0.047989230000000001
$ws.Range("I16").Value = 0.05706583688338699
$ws.Range("C17").Value = 0.00310735
$ws.Range("I17").Value = 0.0064322109649529299
$ws.Range("C18").Value = 0.000029050000000000001
$ws.Range("I18").Value = 0.0065149365542769852
$ws.Range("C19").Value = 0
$ws.Range("I19").Value = 0.027057872087084014
$ws.Range("C20").Value = 0.0031587199999999998
$ws.Range("I20").Value = 0.078546642925202359
$ws.Range("C21").Value = 0.021776719999999999
$ws.Range("I21").Value = 0.011027302924501286
$ws.Range("C22").Value = 0.03625229
$ws.Range("I22").Value = 0.011826997569300394
$ws.Range("C23").Value = 0.043965539999999997
$ws.Range("I23").Value = 0.012378097006121405
$ws.Range("C24").Value = 0.37751600000000002
$ws.Range("I24").Value = 0.11771605876014028
$ws.Range("C25").Value = 0.034526719999999997
$ws.Range("I25").Value = 0.011691046411124525
$ws.Range("C26").Value = 0.02228991
$ws.Range("I26").Value = 0.01121011238942698
$ws.Range("C27").Value = 0.0043419000000000001
$ws.Range("I27").Value = 0.045567523939345392
$ws.Range("I28").Value = 0.025614170375975687
$ws.Range("C29").Value = 0.0011604300000000001
$ws.Range("I29").Value = 0.0035584556609558949
$ws.Range("C30").Value = 0.0094349900000000007
$ws.Range("I30").Value = 0.0039451012011446837
$ws.Range("C31").Value = 0.01376118
$ws.Range("I31").Value = 0.0042561067361226949
$ws.Range("C32").Value = 0.13074809000000001
$ws.Range("I32").Value = 0.03747003333173806
$ws.Range("C33").Value = 0.0068125099999999999
$ws.Range("I33").Value = 0.0035588392309623199
$ws.Range("C34").Value = 0.00068303999999999999
$ws.Range("I34").Value = 0.0035192896763291586
$ws.Range("I35").Value = 0.014846321177627677

# --- load_shapes: refreshed load-shape sample values ---
$ws = $wb.Worksheets.Item("load_shapes")
$ws.Range("B4").Value = 0.04714611872146119
$ws.Range("C4").Value = 0.01861886648600718
$ws.Range("J4").Value = 0.041328001177643389
$ws.Range("O4").Value = 0.29524362667493453
$ws.Range("B5").Value = 0.006735159817351598
$ws.Range("C5").Value = 0.01124563889648543
$ws.Range("J5").Value = 0.0080851488420039929
$ws.Range("O5").Value = 0.047126123672023956
$ws.Range("B6").Value = 0.006735159817351598
$ws.Range("C6").Value = 0.0086053584599192871
$ws.Range("J6").Value = 0.0083943858858496045
$ws.Range("O6").Value = 0.047484009459135335
$ws.Range("B7").Value = 0.006735159817351598
$ws.Range("C7").Value = 0.0085271279284654736
$ws.Range("J7").Value = 0.0085234446018451013
$ws.Range("O7").Value = 0.045197330291383508
$ws.Range("B8").Value = 0.053881278538812784
$ws.Range("C8").Value = 0.079051952034076695
$ws.Range("J8").Value = 0.067674261844725359
$ws.Range("O8").Value = 0.080389650494098852
$ws.Range("B9").Value = 0.006735159817351598
$ws.Range("C9").Value = 0.01016996908899552
$ws.Range("J9").Value = 0.0088413511413772579
$ws.Range("O9").Value = 0.040331772644039976
$ws.Range("B10").Value = 0.006735159817351598
$ws.Range("C10").Value = 0.0091920874458228732
$ws.Range("J10").Value = 0.0085385270462681635
$ws.Range("O10").Value = 0.030789916534887052
$ws.Range("B11").Value = 0.026940639269406392
$ws.Range("C11").Value = 0.01623283527666592
$ws.Range("J11").Value = 0.028553373673719802
$ws.Range("O11").Value = 0.109785485062313
$ws.Range("B12").Value = 0.024771689497716895
$ws.Range("C12").Value = 0.0097827942553597048
$ws.Range("J12").Value = 0.020824568866859339
$ws.Range("O12").Value = 0.27884032566239036
$ws.Range("B13").Value = 0.0035388127853881279
$ws.Range("C13").Value = 0.0059087255218821738
$ws.Range("J13").Value = 0.0039245676009928418
$ws.Range("O13").Value = 0.096611285097571775
$ws.Range("B14").Value = 0.0035388127853881279
$ws.Range("C14").Value = 0.004521459529788099
$ws.Range("J14").Value = 0.0040913868705264866
$ws.Range("O14").Value = 0.096080459952263109
$ws.Range("B15").Value = 0.0035388127853881279
$ws.Range("C15").Value = 0.0044803553522445696
$ws.Range("J15").Value = 0.0041688716615543011
$ws.Range("O15").Value = 0.082359576109102495
$ws.Range("B16").Value = 0.028310502283105023
$ws.Range("C16").Value = 0.041535771407735213
$ws.Range("J16").Value = 0.032744574171450817
$ws.Range("O16").Value = 0.1140433443048483
$ws.Range("B17").Value = 0.0035388127853881279
$ws.Range("C17").Value = 0.0053435430806586613
$ws.Range("J17").Value = 0.0042432288259442244
$ws.Range("O17").Value = 0.083961008243862656
$ws.Range("B18").Value = 0.0035388127853881279
$ws.Range("C18").Value = 0.0048297408613645603
$ws.Range("J18").Value = 0.0041929427472392035
$ws.Range("O18").Value = 0.067301470387042484
$ws.Range("B19").Value = 0.014155251141552512
$ws.Range("C19").Value = 0.0085291168402820937
$ws.Range("J19").Value = 0.014323008290909871
$ws.Range("O19").Value = 0.12565838407205621
$ws.Range("B20").Value = 0.12226027397260274
$ws.Range("C20").Value = 0.048282823260323703
$ws.Range("J20").Value = 0.095989883126765013
$ws.Range("O20").Value = 0.25253238420196467
$ws.Range("B21").Value = 0.017465753424657535
$ws.Range("C21").Value = 0.029162419511224922
$ws.Range("J21").Value = 0.017287412494072642
$ws.Range("O21").Value = 0.13629202385605721
$ws.Range("B22").Value = 0.017465753424657535
$ws.Range("C22").Value = 0.02231559058250256
$ws.Range("J22").Value = 0.018394168048003832
$ws.Range("O22").Value = 0.11899322776779986
$ws.Range("B23").Value = 0.017465753424657535
$ws.Range("C23").Value = 0.022112721577207076
$ws.Range("J23").Value = 0.018615088124530728
$ws.Range("O23").Value = 0.1189849959340632
$ws.Range("B24").Value = 0.13972602739726028
$ws.Range("C24").Value = 0.20499912985108024
$ws.Range("J24").Value = 0.14928254098555391
$ws.Range("O24").Value = 0.14153455122880731
$ws.Range("B25").Value = 0.017465753424657535
$ws.Range("C25").Value = 0.026372970688412108
$ws.Range("J25").Value = 0.018568766960267023
$ws.Range("O25").Value = 0.10810752915343902
$ws.Range("B26").Value = 0.017465753424657535
$ws.Range("C26").Value = 0.023837108122218637
$ws.Range("J26").Value = 0.018425380485267689
$ws.Range("O26").Value = 0.10965799619002992
$ws.Range("B27").Value = 0.069863013698630141
$ws.Range("C27").Value = 0.042095318598811632
$ws.Range("J27").Value = 0.065542910736804488
$ws.Range("O27").Value = 0.23835346524352574
$ws.Range("B28").Value = 0.048744292237442921
$ws.Range("C28").Value = 0.019250014502482003
$ws.Range("J28").Value = 0.037212426865282221
$ws.Range("O28").Value = 0.20134789143260856
$ws.Range("B29").Value = 0.0069634703196347035
$ws.Range("C29").Value = 0.011626846994671376
$ws.Range("J29").Value = 0.0066898936434662487
$ws.Range("O29").Value = 0.085757665349928569
$ws.Range("B30").Value = 0.0069634703196347035
$ws.Range("C30").Value = 0.0088970655263572278
$ws.Range("J30").Value = 0.0071129304901544257
$ws.Range("O30").Value = 0.069062180476435531
$ws.Range("B31").Value = 0.0069634703196347035
$ws.Range("C31").Value = 0.0088161831124812517
$ws.Range("J31").Value = 0.0071906217871427255
$ws.Range("O31").Value = 0.067023222881490874
$ws.Range("B32").Value = 0.055707762557077628
$ws.Range("C32").Value = 0.081731679221672515
$ws.Range("J32").Value = 0.057481011055481458
$ws.Range("O32").Value = 0.091683107910291328
$ws.Range("B33").Value = 0.0069634703196347035
$ws.Range("C33").Value = 0.010514713803876726
$ws.Range("J33").Value = 0.0072010221782435177
$ws.Range("O33").Value = 0.061513546915225481
$ws.Range("B34").Value = 0.0069634703196347035
$ws.Range("C34").Value = 0.009503683630427041
$ws.Range("J34").Value = 0.0071669899959555831
$ws.Range("O34").Value = 0.054986443284197106
$ws.Range("B35").Value = 0.027853881278538814
$ws.Range("C35").Value = 0.016783100879264767
$ws.Range("J35").Value = 0.025371306339946768
$ws.Range("O35").Value = 0.15568071991556032
$ws.Range("J44").Value = 0.046296111212786775
$ws.Range("J45").Value = 0.0067881259361066183
$ws.Range("J46").Value = 0.0068128512704278962
$ws.Range("J47").Value = 0.0068231702796267323
$ws.Range("J48").Value = 0.054544321262162562
$ws.Range("J49").Value = 0.0068485887915071438
$ws.Range("J50").Value = 0.0068243762104431196
$ws.Range("J51").Value = 0.026849692925335147
$ws.Range("J52").Value = 0.024380756632742924
$ws.Range("J53").Value = 0.0035588937149976336
$ws.Range("J54").Value = 0.003572231904236404
$ws.Range("J55").Value = 0.0035784272724857964
$ws.Range("J56").Value = 0.028578932976393311
$ws.Range("J57").Value = 0.0035843725684498889
$ws.Range("J58").Value = 0.0035803518984414517
$ws.Range("J59").Value = 0.014125614503682739
$ws.Range("J60").Value = 0.12051889590551515
$ws.Range("J61").Value = 0.017502793348931428
$ws.Range("J62").Value = 0.017591285014556528
$ws.Range("J63").Value = 0.017608948883938026
$ws.Range("J64").Value = 0.14090052200323819
$ws.Range("J65").Value = 0.017605245232336884
$ws.Range("J66").Value = 0.017593780633893507
$ws.Range("J67").Value = 0.06972279321850966
$ws.Range("J68").Value = 0.048068290744149533
$ws.Range("J69").Value = 0.0069767447410164617
$ws.Range("J70").Value = 0.0070105690439994078
$ws.Range("J71").Value = 0.00701678092362418
$ws.Range("J72").Value = 0.056130732267278459
$ws.Range("J73").Value = 0.0070176124965296805
$ws.Range("J74").Value = 0.00701489142186706
$ws.Range("J75").Value = 0.027795978691315941

$wb.Application.Calculate()
